$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.561.97'
$ws.Range("E2").Value = '  -0.54%  '

$ws.Range("D3").Value = '3.759.85'
$ws.Range("E3").Value = '  -2.00%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '''595.99'
$ws.Range("E5").Value = '  -0.69%  '

$ws.Range("D6").Value = '''170.45'
$ws.Range("E6").Value = '  +1.79%  '

$ws.Range("D7").Value = '3.757.23'
$ws.Range("E7").Value = '  -2.04%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("E9").Value = '  +0.03%  '

$ws.Range("E10").Value = '  +1.02%  '

$ws.Range("D11").Value = '''6.48'
$ws.Range("E11").Value = '  +0.36%  '

$ws.Range("D12").Value = '''0.456'
$ws.Range("E12").Value = '  -0.38%  '

$ws.Range("D13").Value = '''0.0000276'
$ws.Range("E13").Value = '  +5.96%  '

$ws.Range("D14").Value = '''36.78'
$ws.Range("E14").Value = '  -0.37%  '

$ws.Range("D15").Value = '4.384.47'
$ws.Range("E15").Value = '  -2.05%  '

$ws.Range("D16").Value = '3.799.87'
$ws.Range("E16").Value = '  -1.06%  '

$ws.Range("D17").Value = '''18.92'
$ws.Range("E17").Value = '  +3.86%  '

$ws.Range("D18").Value = '67.505.64'
$ws.Range("E18").Value = '  -0.81%  '

$ws.Range("D19").Value = '''7.25'
$ws.Range("E19").Value = '  -1.77%  '

$ws.Range("E20").Value = '  +0.89%  '

$ws.Range("D21").Value = '''10.56'
$ws.Range("E21").Value = '  -3.93%  '

$ws.Range("D22").Value = '''470.11'
$ws.Range("E22").Value = '  +1.04%  '

$ws.Range("D23").Value = '''0.722'
$ws.Range("E23").Value = '  -1.40%  '

$ws.Range("D24").Value = '''0.0000149'
$ws.Range("E24").Value = '  -7.00%  '

$ws.Range("D25").Value = '''83.94'
$ws.Range("E25").Value = '  +1.34%  '

$ws.Range("D26").Value = '''2.24'
$ws.Range("E26").Value = '  +0.61%  '

$ws.Range("E27").Value = '  +0.50%  '

$ws.Range("D28").Value = '''10.43'
$ws.Range("E28").Value = '  +3.81%  '

$ws.Range("E30").Value = '  -1.90%  '

$ws.Range("D31").Value = '3.906.68'
$ws.Range("E31").Value = '  -1.91%  '

$ws.Range("D32").Value = '''7.72'
$ws.Range("E32").Value = '  +1.31%  '

$ws.Range("E33").Value = '  -2.49%  '

$ws.Range("D34").Value = '''30.47'
$ws.Range("E34").Value = '  -2.35%  '

$ws.Range("D35").Value = '''9.17'
$ws.Range("E35").Value = '  -4.05%  '

$ws.Range("D36").Value = '3.719.56'
$ws.Range("E36").Value = '  -1.99%  '

$ws.Range("D37").Value = '''3.86'
$ws.Range("E37").Value = '  +7.60%  '

$ws.Range("D38").Value = '''0.106'
$ws.Range("E38").Value = '  +1.36%  '

$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").Value = '''5.90'
$ws.Range("E39").Value = '  -0.21%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '''0.138'
$ws.Range("E40").Value = '  -1.65%  '

$ws.Range("E41").Value = '  -1.87%  '

$ws.Range("D42").Value = '''0.999'
$ws.Range("E42").Value = '  +0.02%  '

$ws.Range("D43").Value = '''0.315'
$ws.Range("E43").Value = '  +0.45%  '

$ws.Range("D45").Value = '''8.75'
$ws.Range("E45").Value = '  +0.87%  '

$ws.Range("D46").Value = '''1.96'
$ws.Range("E46").Value = '  -0.81%  '

$ws.Range("D47").Value = '''45.85'
$ws.Range("E47").Value = '  -2.23%  '

$ws.Range("D48").Value = '''398.48'
$ws.Range("E48").Value = '  -5.43%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = '''141.85'
$ws.Range("E49").Value = '  -0.21%  '

$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").Value = '''0.000270'
$ws.Range("E50").Value = '  -6.87%  '

$ws.Range("E51").Value = '  -0.44%  '
